# Apply weekly crypto price/volume refresh (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.779.52"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3
$ws.Range("D3").Value = "2.619.19"
$ws.Range("E3").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'594.44"
$ws.Range("E5").Value = "  -0.16%  "

# Row 6
$ws.Range("D6").Value = "'150.83"
$ws.Range("E6").Value = "  +0.60%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("E9").Value = "  +4.51%  "

# Row 10
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'5.79"
$ws.Range("E10").Value = "  +2.15%  "

# Row 11
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.394"
$ws.Range("E11").Value = "  +3.07%  "

# Row 12
$ws.Range("E12").Value = "  +1.09%  "

# Row 13
$ws.Range("E13").Value = "  +0.79%  "

# Row 14
$ws.Range("D14").Value = "3.091.05"
$ws.Range("E14").Value = "  -0.04%  "

# Row 15
$ws.Range("D15").Value = "63.706.93"
$ws.Range("E15").Value = "  +0.26%  "

# Row 16
$ws.Range("D16").Value = "'0.0000169"
$ws.Range("E16").Value = "  +13.56%  "

# Row 17
$ws.Range("D17").Value = "2.630.67"
$ws.Range("E17").Value = "  -0.54%  "

# Row 18
$ws.Range("E18").Value = "  -0.35%  "

# Row 19
$ws.Range("E19").Value = "  +3.07%  "

# Row 20
$ws.Range("D20").Value = "'347.81"
$ws.Range("E20").Value = "  -0.07%  "

# Row 21
$ws.Range("E21").Value = "  +2.29%  "

# Row 23
$ws.Range("E23").Value = "  +1.76%  "

# Row 24
$ws.Range("E24").Value = "  -3.17%  "

# Row 25
$ws.Range("E25").Value = "  +0.49%  "

# Row 26
$ws.Range("E26").Value = "  +0.02%  "

# Row 27
$ws.Range("E27").Value = "  +0.86%  "

# Row 28
$ws.Range("D28").Value = "'546.87"
$ws.Range("E28").Value = "  -3.04%  "

# Row 29
$ws.Range("E29").Value = "  -1.52%  "

# Row 30
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.23%  "

# Row 31
$ws.Range("E31").Value = "  +2.09%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0904"
$ws.Range("E32").Value = "  +7.43%  "

# Row 33
$ws.Range("D33").Value = "'1.83"
$ws.Range("E33").Value = "  +4.66%  "

# Row 34
$ws.Range("E34").Value = "  +4.53%  "

# Row 35
$ws.Range("E35").Value = "  +1.13%  "

# Row 36
$ws.Range("D36").Value = "'0.421"
$ws.Range("E36").Value = "  +3.07%  "

# Row 37
$ws.Range("D37").Value = "'164.94"
$ws.Range("E37").Value = "  -2.24%  "

# Row 38
$ws.Range("D38").Value = "'19.91"
$ws.Range("E38").Value = "  +3.19%  "

# Row 39
$ws.Range("E39").Value = "  +2.07%  "

# Row 40
$ws.Range("E40").Value = "  -0.09%  "

# Row 41
$ws.Range("E41").Value = "  -0.02%  "

# Row 42
$ws.Range("D42").Value = "'168.22"
$ws.Range("E42").Value = "  -0.54%  "

# Row 43
$ws.Range("E43").Value = "  +4.72%  "

# Row 44
$ws.Range("D44").Value = "'23.16"
$ws.Range("E44").Value = "  +8.82%  "

# Row 45
$ws.Range("E45").Value = "  -2.03%  "

# Row 46
$ws.Range("E46").Value = "  +10.95%  "

# Row 47
$ws.Range("E47").Value = "  +1.23%  "

# Row 48
$ws.Range("E48").Value = "  +1.49%  "

# Row 49
$ws.Range("D49").Value = "'0.0971"
$ws.Range("E49").Value = "  +0.40%  "

# Row 50
$ws.Range("D50").Value = "'19.23"
$ws.Range("E50").Value = "  +0.62%  "

# Row 51
$ws.Range("E51").Value = "  +19.20%  "
